# S-01004 - Avance en los datos de los graficos de performance de agentes
#
# "Nico" logged two more days of work (4h each) on story S-01004
# ("Desarrollo Metricas Agentes") in the "Horas insumidas" sheet, and two
# pre-existing rows had their dates corrected. This ripples through the
# SUMIF-based totals on "Earned Value" (and, from there, "Estadísticas").

$wb = $excel.ActiveWorkbook

$hs = $wb.Worksheets.Item("Horas insumidas")
$ev = $wb.Worksheets.Item("Earned Value")

# --- Correct the dates on the two existing rows (17-Oct -> 14-Oct) ---
$hs.Range("B68").Value = 40465
$hs.Range("B69").Value = 40465

# --- Log the new work on S-01004 : "Desarrollo Metricas Agentes" ---
$hs.Range("B70").Value = 40465
$hs.Range("B70").NumberFormat = "d-mmm"
$hs.Range("C70").Value = "Nico"
$hs.Range("D70").Value = "Desarrollo Metricas Agentes"
$hs.Range("E70").Value = "S-01004"
$hs.Range("F70").Value = 4

$hs.Range("B71").Value = 40466
$hs.Range("B71").NumberFormat = "d-mmm"
$hs.Range("C71").Value = "Nico"
$hs.Range("D71").Value = "Desarrollo Metricas Agentes"
$hs.Range("E71").Value = "S-01004"
$hs.Range("F71").Value = 4

$wb.Application.CalculateFull() | Out-Null

# The new rows sit inside the SUMIF scan range ('Horas insumidas'!$E$6:$E$131)
# that "Earned Value"!G2:G18 reads, so G4 already refreshes to 16. The
# dependent SUM(...) totals a row below (G21/H21) need their formulas
# re-seated (clear, then re-enter) so they pick up that new precedent.
$g21Formula = $ev.Range("G21").Formula()
$h21Formula = $ev.Range("H21").Formula()
$ev.Range("G21:H21").ClearContents()
$ev.Range("G21").Formula = $g21Formula
$ev.Range("H21").Formula = $h21Formula

# --- Reflect where the user was last working in "Horas insumidas" ---
$hs.Activate()
$hs.Range("D73").Select() | Out-Null

$wb.Application.CalculateFull() | Out-Null
